# This edit re-distributes the observation records currently stored in
# rows 8-24 of the "Artfynd" sheet: the full content of each row is
# replaced by the full content that another row in that same block used
# to hold (a permutation of the 17 rows). Row numbers/formatting stay
# fixed; only the cell contents move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row number -> row number whose original content should end
# up there.
$rowMap = @{
    8  = 22
    9  = 11
    10 = 19
    11 = 20
    12 = 21
    13 = 10
    14 = 8
    15 = 16
    16 = 18
    17 = 13
    18 = 9
    19 = 24
    20 = 12
    21 = 15
    22 = 23
    23 = 17
    24 = 14
}

# The row data is copied in two column chunks, A:X and AC:AY, skipping
# Y:AB (Startdatum/Starttid/Slutdatum/Sluttid). Those four columns hold
# the exact same literal text ("2022-05-30"/"00:00") on every row in this
# block, so they never actually change - and skipping them avoids Excel's
# automatic text->date/time conversion that setting .Value2 would trigger.
$chunks = @(
    @{ First = "A";  Last = "X" },
    @{ First = "AC"; Last = "AY" }
)

# 1) Snapshot every affected row's current values before any writes,
#    otherwise a later write could clobber data still needed as a source.
$snapshots = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshots.ContainsKey($srcRow)) {
        $rowData = @{}
        foreach ($chunk in $chunks) {
            $rng = $ws.Range("$($chunk.First)$srcRow`:$($chunk.Last)$srcRow")
            $rowData[$chunk.First] = $rng.Value2
        }
        $snapshots[$srcRow] = $rowData
    }
}

# 2) Write each target row using the snapshot captured in step 1.
foreach ($targetRow in $rowMap.Keys) {
    $srcRow = $rowMap[$targetRow]
    $rowData = $snapshots[$srcRow]
    foreach ($chunk in $chunks) {
        $destRng = $ws.Range("$($chunk.First)$targetRow`:$($chunk.Last)$targetRow")
        $destRng.Value2 = $rowData[$chunk.First]
    }
}

Write-Output "Row contents permuted for rows 8-24."
